# Update the ripspillningshögar distance-sampling summary sheet:
# - keep only Label/Estimate columns (rename + new values)
# - drop se/cv/lcl/ucl/df columns and the Total row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels
$ws.Range("A1").Value = "lya"
$ws.Range("B1").Value = "uppskattat_antal_ripspillningshögar"

# New Estimate values for rows 2-11
$values = @(
    36724.77560238716,
    21970.30615143725,
    19807.1803297371,
    15764.00497871487,
    20001.64072568122,
    6078.930024471744,
    17081.29097785449,
    17665.26673778114,
    8540.645488927245,
    11956.90368449814
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Remove the now-unused se/cv/lcl/ucl/df columns (C:G)
$ws.Range("C1:G12").Delete()

# Remove the Total row (row 12)
$ws.Range("A12:B12").Delete()
